$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "(Clock Cycles)" header above the BEAGLE BONE BLACK memset table (H4)
# was inconsistent with its memmove counterpart (B4, "(values in Clock
# Cycles)"); unify both to the fuller label.
$ws.Range("H4").Value = "(values in Clock Cycles)"

# Likewise the microsecond-unit labels for the RTC-timed (KL25Z) tables were
# just "uS" -- make them consistent with the other unit labels used
# elsewhere in the sheet ("(values in uS)").
$ws.Range("B12").Value = "(values in uS)"
$ws.Range("H12").Value = "(values in uS)"

# Move the active cell/selection to where the author left off editing.
[void]$ws.Range("I27").Select()
